$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.536.71'
$ws.Range("E2").Value = '  -0.66%  '

$ws.Range("D3").Value = '2.282.90'
$ws.Range("E3").Value = '  -1.23%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '95.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.71'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.620'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  -2.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0930'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.89'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").Value = '2.625.22'
$ws.Range("E14").Value = '  -1.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.839'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.90%  '

$ws.Range("D17").Value = '2.281.53'
$ws.Range("E17").Value = '  -1.32%  '

$ws.Range("D18").Value = '43.508.55'
$ws.Range("E18").Value = '  -0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000109'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.34%  '

$ws.Range("E20").Value = '  -1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +13.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.59%  '

$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.17%  '

$ws.Range("E30").Value = '  -2.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.94%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0892'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.09%  '

$ws.Range("E34").Value = '  -4.13%  '

$ws.Range("E35").Value = '  -0.87%  '

$ws.Range("E36").Value = '  -3.31%  '

$ws.Range("E37").Value = '  -1.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.16%  '

$ws.Range("E39").Value = '  -3.24%  '

$ws.Range("E40").Value = '  +1.10%  '

$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.01%  '

$ws.Range("E44").Value = '  +1.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.74'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.101'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.00%  '

$ws.Range("E48").Value = '  -1.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.23%  '

$ws.Range("E50").Value = '  +8.06%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.423'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.54%  '
